$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

for ($row = 2; $row -le 40; $row++) {
    $ws.Cells.Item($row, 30).Value = 103
    $ws.Cells.Item($row, 31).Value = 58
    $ws.Cells.Item($row, 32).Value = 0
}
